# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.002.97"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.429.43"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'410.82"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'130.28"
$ws.Range("D7").Value = "'0.636"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").Value = "'43.63"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "'0.0000226"
$ws.Range("E12").Value = "  +15.48%  "
$ws.Range("E13").Value = "  +4.80%  "
$ws.Range("D14").Value = "3.969.37"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'21.24"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").Value = "3.426.50"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'12.36"
$ws.Range("E18").Value = "  +6.88%  "
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "61.958.12"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "'518.81"
$ws.Range("E21").Value = "  +31.08%  "
$ws.Range("D22").Value = "'92.75"
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("D23").Value = "'3.32"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("D24").Value = "'13.43"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "'3.34"
$ws.Range("D26").Value = "'34.94"
$ws.Range("E26").Value = "  +9.15%  "
$ws.Range("D27").Value = "'9.28"
$ws.Range("E27").Value = "  +9.12%  "
$ws.Range("D28").Value = "'7.66"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'12.17"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "'0.167"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").Value = "'42.04"
$ws.Range("D34").Value = "'59.28"
$ws.Range("E34").Value = "  +13.18%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.140"
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = "  +18.34%  "
$ws.Range("D41").Value = "'148.20"
$ws.Range("E41").Value = "  +5.46%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.95"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'2.11"
$ws.Range("E43").Value = "  +7.04%  "
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "'4.32"
$ws.Range("E45").Value = "  +7.58%  "
$ws.Range("E46").Value = "  +22.38%  "
$ws.Range("D47").Value = "'16.71"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'23.25"
$ws.Range("E48").Value = "  +4.42%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'121.52"
$ws.Range("E49").Value = "  +26.32%  "
$ws.Range("E50").Value = "  +19.97%  "
$ws.Range("D51").Value = "2.141.32"
$ws.Range("E51").Value = "  +0.75%  "
